# Add blurbs for each path
$wb = $excel.ActiveWorkbook

# --- Rename the original sheet, add a new "Paths" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Projects"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Paths"

# --- Move the selection on the Projects sheet from D4:M4 to D1:M1 ---
$ws1.Range("D1:M1").Select() | Out-Null

# --- Populate the Paths sheet ---
$pathNames = @(
    "Dynamic Leadership",
    "Effective Coaching",
    "Innovative Planning",
    "Leadership Development",
    "Motivational Strategies",
    "Persuasive Influence",
    "Presentation Mastery",
    "Strategic Relationships",
    "Team Collaboration",
    "Visionary Communication"
)

$blurbs = @(
    "This path helps you build your skills as a strategic leader. The projects on this path focus on understanding leadership and communication styles, the effect of conflict on a group and the skills needed to defuse and direct conflict. These projects also emphasize the development of strategies to facilitate change in an organization or group, interpersonal communication and public speaking. This path culminates in a project focused on applying your leadership skills.",
    "This path helps you build your skills as a positive communicator and leader. The projects on this path focus on understanding and building consensus, contributing to the development of others by coaching and establishing strong public speaking skills. Each project emphasizes the importance of effective interpersonal communication. This path culminates in a “High Performance Leadership” project of your design.",
    "This path helps you build your skills as a public speaker and leader. The projects on this path focus on developing a strong connection with audience members when you present, speech writing and speech delivery. The projects contribute to building an understanding of the steps to manage a project, as well as creating innovative solutions. This path culminates in a “High Performance Leadership” project of your design.",
    "This path helps you build your skills as an effective communicator and leader. The projects on this path focus on learning how to manage time, as well as how to develop and implement a plan. Public speaking and leading a team are emphasized in all projects. This path culminates in the planning and execution of an event that will allow you to apply everything you learned.",
    "This path helps you build your skills as a powerful and effective communicator. The projects focus on learning strategies for building connections with the people around you, understanding motivation and successfully leading small groups to accomplish tasks. This path culminates in a comprehensive team-building project that brings all of your skills together—including public speaking.",
    "This path helps you build your skills as an innovative communicator and leader. The projects on this path focus on how to negotiate a positive outcome together with building strong interpersonal communication and public speaking skills. Each project emphasizes developing leadership skills to use in complex situations, as well as creating innovative solutions to challenges. This path culminates in a “High Performance Leadership” project of your design.",
    "This path helps you build your skills as an accomplished public speaker. The projects on this path focus on learning how an audience responds to you and improving your connection with audience members. The projects contribute to developing an understanding of effective public speaking technique, including speech writing and speech delivery. This path culminates in an extended speech that will allow you to apply what you learned.",
    "This path helps you build your skills as a leader in communication. The projects on this path focus on understanding diversity, building personal and/or professional connections with a variety of people and developing a public relations strategy. Communicating well interpersonally and as a public speaker is emphasized in each project. The path culminates in a project to apply your skills as a leader in a volunteer organization.",
    "This path helps you build your skills as a collaborative leader. The projects on this path focus on active listening, motivating others and collaborating with a team. Each project contributes to building interpersonal communication and public speaking skills. This path culminates in a project focused on applying your leadership skills.",
    "This path helps you build your skills as a strategic communicator and leader. The projects on this path focus on developing your skills for sharing information with a group, planning communications and creating innovative solutions. Speech writing and speech delivery are emphasized in each project. This path culminates in the development and launch of a long-term personal or professional vision."
)

$ws2.Range("A1").Value = "Pathname"
$ws2.Range("B1").Value = "Blurb"

for ($i = 0; $i -lt $pathNames.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $pathNames[$i]
    $ws2.Cells.Item($row, 2).Value = $blurbs[$i]
}

$ws2.Columns.Item(1).ColumnWidth = 22
$ws2.Columns.Item(2).ColumnWidth = 145.1640625

$ws2.Range("A1:A11").VerticalAlignment = -4160
$ws2.Range("B2:B11").WrapText = $true
$ws2.Range("B2:B11").VerticalAlignment = -4160

$heights = @(68, 51, 51, 51, 51, 68, 51, 51, 51, 51)
for ($i = 0; $i -lt $heights.Length; $i++) {
    $ws2.Rows.Item($i + 2).RowHeight = $heights[$i]
}

$ws2.Range("B1:B2").Font.Size = 13
$ws2.Range("B1:B2").Font.Color = 3355443
$ws2.Range("B1:B2").Font.Name = "Arial"
$ws2.Range("B1").WrapText = $true

$ws2.Range("B2").Select() | Out-Null
$ws2.Range("A1").Select() | Out-Null
